$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "72.237.93"
$ws.Range("E2").Value = "  +4.21%  "
$ws.Range("D3").Value = "4.033.05"
$ws.Range("E3").Value = "  +3.49%  "
$ws.Range("E4").Value = "  +0.04%  "
$ws.Range("D5").Value = "520.41"
$ws.Range("E5").Value = "  -2.04%  "
$ws.Range("D6").Value = "147.02"
$ws.Range("E6").Value = "  +1.64%  "
$ws.Range("D7").Value = "0.733"
$ws.Range("E7").Value = "  +19.78%  "
$ws.Range("D8").Value = "4.024.65"
$ws.Range("E8").Value = "  +3.44%  "
$ws.Range("E9").Value = "  +0.14%  "
$ws.Range("D10").Value = "0.777"
$ws.Range("E10").Value = "  +8.13%  "
$ws.Range("E11").Value = "  +1.72%  "
$ws.Range("E12").Value = "  -1.79%  "
$ws.Range("D13").Value = "48.29"
$ws.Range("E13").Value = "  +14.66%  "
$ws.Range("D14").Value = "11.11"
$ws.Range("E14").Value = "  +8.31%  "
$ws.Range("D15").Value = "4.678.34"
$ws.Range("E15").Value = "  +3.58%  "
$ws.Range("D16").Value = "4.053.56"
$ws.Range("E16").Value = "  +3.70%  "
$ws.Range("D17").Value = "21.23"
$ws.Range("E17").Value = "  +7.34%  "
$ws.Range("E18").Value = "  +1.37%  "
$ws.Range("D19").Value = "1.21"
$ws.Range("E19").Value = "  +0.35%  "
$ws.Range("E20").Value = "  -0.51%  "
$ws.Range("D21").Value = "72.168.98"
$ws.Range("E21").Value = "  +4.18%  "
$ws.Range("D22").Value = "444.54"
$ws.Range("D23").Value = "105.22"
$ws.Range("E23").Value = "  +19.70%  "
$ws.Range("E24").Value = "  +5.36%  "
$ws.Range("D25").Value = "14.98"
$ws.Range("E25").Value = "  +6.04%  "
$ws.Range("D26").Value = "4.01"
$ws.Range("E26").Value = "  -0.54%  "
$ws.Range("D27").Value = "11.50"
$ws.Range("E27").Value = "  -0.23%  "
$ws.Range("D28").Value = "11.03"
$ws.Range("E28").Value = "  +4.49%  "
$ws.Range("D29").Value = "37.79"
$ws.Range("E29").Value = "  +4.01%  "
$ws.Range("D30").Value = "5.81"
$ws.Range("E30").Value = "  +2.30%  "
$ws.Range("D31").Value = "3.25"
$ws.Range("E31").Value = "  +14.79%  "
$ws.Range("D32").Value = "13.70"
$ws.Range("E32").Value = "  +3.95%  "
$ws.Range("E33").Value = "  +2.94%  "
$ws.Range("D34").Value = "675.56"
$ws.Range("E34").Value = "  -1.82%  "
$ws.Range("E35").Value = "  +14.21%  "
$ws.Range("D36").Value = "66.76"
$ws.Range("E36").Value = "  -2.55%  "
$ws.Range("D37").Value = "42.37"
$ws.Range("E37").Value = "  +6.28%  "
$ws.Range("E38").Value = "  -0.19%  "
$ws.Range("D39").Value = "0.425"
$ws.Range("E40").Value = "  +6.87%  "
$ws.Range("E41").Value = "  +1.59%  "
$ws.Range("E42").Value = "  +0.06%  "
$ws.Range("D43").Value = "0.0501"
$ws.Range("E43").Value = "  +3.51%  "
$ws.Range("D44").Value = "0.998"
$ws.Range("E44").Value = "  -0.20%  "
$ws.Range("B45").Value = "Stellar"
$ws.Range("C45").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D45").Value = "0.162"
$ws.Range("E45").Value = "  +15.10%  "
$ws.Range("B46").Value = "WEMIXToken"
$ws.Range("C46").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("D46").Value = "3.26"
$ws.Range("E46").Value = "  +1.56%  "
$ws.Range("B47").Value = "ApeXProtocol"
$ws.Range("C47").Value = "https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex"
$ws.Range("D47").Value = "3.53"
$ws.Range("E47").Value = "  +4.08%  "
$ws.Range("B48").Value = "Fetch.AI"
$ws.Range("C48").Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
$ws.Range("D48").Value = "2.71"
$ws.Range("E48").Value = "  -2.46%  "
$ws.Range("D49").Value = "3.06"
$ws.Range("E49").Value = "  +2.32%  "
$ws.Range("D50").Value = "9.27"
$ws.Range("E50").Value = "  +8.59%  "
$ws.Range("E51").Value = "  +2.73%  "
